# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet right before the "总计" sheet (mirrors the
#    same per-fund holdings layout used by the other quarterly sheets).
# 2) Insert a new first data row on "总计" summarising the new quarter and
#    shift the previously-existing summary rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the "2022-Q1" worksheet, positioned immediately before "总计".
# ---------------------------------------------------------------------------
$zongjiBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zongjiBefore)
$q1.Name = "2022-Q1"

# NOTE: worksheet variables captured before a sheet is inserted can end up
# referring to whatever sheet now sits at that same *position* rather than
# the original sheet, so re-fetch "总计" by name now that the tab order has
# changed.
$zongji = $wb.Worksheets.Item("总计")

# Header row (matches the other quarterly sheets, e.g. "2021-Q4").
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Carry over the header formatting (bold + border + centered) from another
# quarterly sheet so the new tab matches its siblings.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A15").PasteSpecial(-4122)

$fundRows = @(
    @(0,  "166027", "中欧创业板两年定期开放混合A",      "21.11", "99.64", "4.77", "1.0069", 7),
    @(1,  "233007", "大摩卓越成长混合",                  "4.96",  "92.50", "6.55", "0.3249", 5),
    @(2,  "010314", "摩根士丹利华鑫内需增长混合",        "4.94",  "94.09", "5.75", "0.2840", 8),
    @(3,  "009791", "中欧创业板两年定期开放混合C",      "5.21",  "99.64", "4.77", "0.2485", 7),
    @(4,  "163302", "大摩资源优选混合(LOF)",            "5.82",  "81.78", "4.15", "0.2415", 4),
    @(5,  "000586", "景顺长城中小板创业板精选股票",      "2.42",  "94.15", "6.64", "0.1607", 4),
    @(6,  "009499", "景顺长城安鑫回报一年持有期混合A",  "1.11",  "29.52", "4.95", "0.0549", 4),
    @(7,  "582003", "东吴配置优化灵活配置混合",          "1.04",  "90.74", "3.15", "0.0328", 9),
    @(8,  "004694", "天弘策略精选灵活配置混合A",        "1.11",  "80.93", "2.31", "0.0256", 10),
    @(9,  "350001", "天治财富增长混合",                  "0.98",  "69.00", "2.35", "0.0230", 10),
    @(10, "009755", "景顺长城安鑫回报一年持有期混合C",  "0.21",  "29.52", "4.95", "0.0104", 4),
    @(11, "005104", "富荣福康混合A",                      "0.08",  "87.88", "3.03", "0.0024", 7),
    @(12, "004748", "天弘策略精选灵活配置混合C",        "0.08",  "80.93", "2.31", "0.0018", 10),
    @(13, "005105", "富荣福康混合C",                      "0.04",  "87.88", "3.03", "0.0012", 7)
)

$row = 2
foreach ($fund in $fundRows) {
    $q1.Range("A$row").Value = $fund[0]
    $q1.Range("B$row").Value = "'" + $fund[1]
    $q1.Range("C$row").Value = $fund[2]
    $q1.Range("D$row").Value = "'" + $fund[3]
    $q1.Range("E$row").Value = "'" + $fund[4]
    $q1.Range("F$row").Value = "'" + $fund[5]
    $q1.Range("G$row").Value = "'" + $fund[6]
    $q1.Range("H$row").Value = $fund[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2. Shift "总计" rows down by one and insert the 2022-Q1 summary at the top.
#    (Literal values are used instead of reading `.Value` back from another
#    range - the read side of that property isn't reliable in this host.)
# ---------------------------------------------------------------------------
# Row 7 is brand new (the sheet previously only went down to row 6) so it
# needs the index-column style (bold/centered/bordered) copied over first.
$zongji.Range("A6").Copy()
$zongji.Range("A7").PasteSpecial(-4122)

$zongji.Range("A7").Value = 5
$zongji.Range("B7").Value = "2020-Q4"
$zongji.Range("C7").Value = 29
$zongji.Range("D7").Value = 6.08

$zongji.Range("A6").Value = 4
$zongji.Range("B6").Value = "2021-Q1"
$zongji.Range("C6").Value = 33
$zongji.Range("D6").Value = 5.13

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2021-Q2"
$zongji.Range("C5").Value = 26
$zongji.Range("D5").Value = 4.48

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q3"
$zongji.Range("C4").Value = 9
$zongji.Range("D4").Value = 1.95

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2021-Q4"
$zongji.Range("C3").Value = 24
$zongji.Range("D3").Value = 3.02

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 14
$zongji.Range("D2").Value = 2.42
